# add x64 config file
$wb = $excel.ActiveWorkbook

# Rename "Sheet2" to "bugs"
$sheetBugs = $wb.Worksheets.Item("Sheet2")
$sheetBugs.Name = "bugs"

$sheetPlans = $wb.Worksheets.Item("plans")

# Add new rows of data to the "plans" sheet
$sheetPlans.Range("A2").Value = "custom menu"
$sheetPlans.Range("B1").Copy() | Out-Null
$sheetPlans.Range("B2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$sheetPlans.Range("B2").Value = 42423

$sheetPlans.Range("A3").Value = "custom toolbar"
$sheetPlans.Range("A4").Value = "custom others"

# Update selections on each sheet
$sheetPlans.Range("H8").Select() | Out-Null

# Make "bugs" the active sheet/tab
$sheetBugs.Activate()
$sheetBugs.Range("H10").Select() | Out-Null
